$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# --- New shared strings must be created in this exact order so sharedStrings.xml
# appends them as: 20513841, 1989636238, USUARIOS41, MORAM12, msg1, msg2 ---

# 1) F5 = "20513841" (numeric-looking text, keep default @ style, no quote-prefix)
$ws.Range("F5").Value2 = "20513841"

# 2) F7 = "1989636238" (numeric-looking text, keep default @ style, no quote-prefix)
$ws.Range("F7").Value2 = "1989636238"

# 3) G5 = "USUARIOS41" typed with a leading apostrophe (quote-prefixed text) -> new style
$ws.Range("G5").Value2 = "'USUARIOS41"

# 4) G6 = "USUARIOS41" (reuses the shared string, plain style)
$ws.Range("G6").Value2 = "USUARIOS41"

# 5) G8 = "MORAM12"
$ws.Range("G8").Value2 = "MORAM12"

# 6/7) B5 and B7 get a yellow highlight fill -> new style (quote-prefix + yellow fill)
$ws.Range("B5").Interior.Color = 65535
$ws.Range("B7").Interior.Color = 65535

# 8/9/10) E5, E6, E8 = first error message
$msg1 = "Usuario o clave inválida. Inténtalo nuevamente"
$ws.Range("E5").Value2 = $msg1
$ws.Range("E6").Value2 = $msg1
$ws.Range("E8").Value2 = $msg1

# 11/12) E7 gets word-wrap (new style) then the long message
$ws.Range("E7").WrapText = $true
$ws.Range("E7").Value2 = "La clave que usas en el cajero está bloqueada. Debes activarla en la Sucursal Física. Para mayor información comunícate con la Sucursal Telefónica."

# Numeric transaction codes
$ws.Range("H5").Value2 = 1234
$ws.Range("H6").Value2 = 4567
$ws.Range("H7").Value2 = 1234
$ws.Range("H8").Value2 = 1234

# Column E is now much wider to fit the long messages
$ws.Columns.Item(5).ColumnWidth = 138.42

# Update the remembered selection/scroll position
$ws.Range("E12").Select() | Out-Null
